# Chapter_4_Table_S4.4.xlsx - final supplementary upload edit
# - Shorten the "Study" column citations to compact codes
#   (also fixes a stray duplicate/typo string "Wiarda et al., (2020" that
#    was missing its closing parenthesis, unifying it with the correct one)
# - Restyle the table
# - Tidy the header row height (auto-fit to match the rest of the table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$studyCol = $ws.Range("D1:D70")

# NOTE: order matters - replace the full/correct strings first so that the
# shorter, truncated typo variant below can't accidentally partial-match them.
$studyCol.Replace("O'Grady et al., (2025)", "OGR25-BTB", `
    [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole, `
    [Microsoft.Office.Interop.Excel.XlSearchOrder]::xlByRows, $true, $false, $false)

$studyCol.Replace("Wiarda et al., (2020)", "WIA20-BTB", `
    [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole, `
    [Microsoft.Office.Interop.Excel.XlSearchOrder]::xlByRows, $true, $false, $false)

# Fix/unify the stray duplicate string (missing closing paren) onto the same code
$studyCol.Replace("Wiarda et al., (2020", "WIA20-BTB", `
    [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole, `
    [Microsoft.Office.Interop.Excel.XlSearchOrder]::xlByRows, $true, $false, $false)

$studyCol.Replace("McLoughlin et al., (2021)", "MCL21-BTB", `
    [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole, `
    [Microsoft.Office.Interop.Excel.XlSearchOrder]::xlByRows, $true, $false, $false)

$studyCol.Replace("McLoughlin et al., (2014)", "MCL14-BTB", `
    [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole, `
    [Microsoft.Office.Interop.Excel.XlSearchOrder]::xlByRows, $true, $false, $false)

# Update the table's visual style
$tbl = $ws.ListObjects.Item(1)
$tbl.TableStyle = "TableStyleMedium15"

# Header row (row 2) height reverts to the same auto-fit height as the data rows
$ws.Rows.Item(2).RowHeight = 24.95
